$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: merge "El knn" + " es un mé" + "todo de clasificación supervisado" +
# " que sirve para estimar la función de densidad F(x/" into a single run
# (pure run-consolidation, the concatenated text is unchanged).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("El knn es un mé", $true, $false, $false, $false, $false, $true, 1, $false, "El knn es un mé", 2)

# ---------------------------------------------------------------------------
# Hunk 2: merge "C" + "j" -> "Cj" (first split occurrence, right after
# "x por cada clase ") and merge ". " + "Este es un método..." into a single
# run, without disturbing the surrounding w:proofErr spellStart/spellEnd
# markers that wrap "Cj".
# ---------------------------------------------------------------------------
$probe1 = $d.Content.Duplicate
$probe1.Find.Execute("cada clase Cj", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$cjRange1 = $d.Range($probe1.End - 2, $probe1.End)
$cjRange1.Find.Execute("Cj", $true, $false, $false, $false, $false, $true, 1, $false, "Cj", 2)

# ---------------------------------------------------------------------------
# Hunk 3: merge "C" + "j" -> "Cj" (second split occurrence, right after
# "pertenezca a la clase ").
# ---------------------------------------------------------------------------
$probe2 = $d.Content.Duplicate
$probe2.Find.Execute("pertenezca a la clase Cj", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$cjRange2 = $d.Range($probe2.End - 2, $probe2.End)
$cjRange2.Find.Execute("Cj", $true, $false, $false, $false, $false, $true, 1, $false, "Cj", 2)

# ---------------------------------------------------------------------------
# Hunk 4: insert nine new paragraphs right after the paragraph holding the
# two inline images (precision/recall curve + ROC curve) and right before
# the empty "_GoBack" bookmark paragraph that precedes the "EDA:" heading.
# ---------------------------------------------------------------------------
$lastShape = $d.InlineShapes.Item($d.InlineShapes.Count)
$imgPara = $lastShape.Range.Paragraphs.Item(1)
$insertionPoint = $imgPara.Range.Duplicate
$insertionPoint.Collapse(0)
$newParagraphsXml = '<w:p><w:r><w:t>Observando estos casos podemos ver la compensación entre la precisión y exhaustividad para diferentes umbrales, de manera que para el primero se observa hasta un 0.7 de exhaustividad aproximadamente un alto grado de precisión y exhaustividad relacionado así con la tasa baja de falsos positivos y una tasa baja de falsos negativos. A partir de este valor hacia arriba se observa sobretodo un incremento de falsos negativos.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Para la curva ROC podemos observar una correcta proporción de verdaderos y falsos positivos indicándonos que nuestro modelo tiene una alta precisión para los valores de entrada que ha </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>obtenido siendo esta curva proporcional al número de datos registrado, ya que a medida que aumentan los falsos negativos siguen estando los valores de verdaderos positivos muy altos para ambos modelos (Regresor logístico y svm).</w:t></w:r></w:p><w:p><w:r><w:t>Con este análisis de ambas gráficas, se procede a la realización del modelo svm con varios tipos de funciones kernel y variaciones entre las variables slack.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Para explicar la afectación del kernel en las SVM </w:t></w:r><w:r><w:t>Se tiene que tener en cuenta que en muchos casos los grupos a clasificar no serán linealmente separables en el espacio original, por lo tanto una solución que ofrece</w:t></w:r><w:r><w:t xml:space="preserve"> la Maquina de vectores de soporte es</w:t></w:r><w:r><w:t xml:space="preserve"> aumentar la dimensión de los datos, la cual se puede transformar combinando o modificando cualquiera de sus dimensiones. Para hacerlo se utiliza el kernel que se una función que devuelve el resultado del producto entre dos vectores realizado en un nuevo espacio dimensional diferente al espacio original en el que se encontraban.</w:t></w:r></w:p><w:p><w:r><w:t>De forma que nos permite operar en el espacio de características original sin calcular las coordenadas de los datos en un espacio de mayor dimensió</w:t></w:r><w:r><w:t xml:space="preserve">n, ofreciéndonos </w:t></w:r><w:r><w:t xml:space="preserve">en esencia una forma más eficiente y menos costosa de transformar los datos en dimensiones </w:t></w:r><w:r><w:t>más</w:t></w:r><w:r><w:t xml:space="preserve"> altas.</w:t></w:r></w:p><w:p><w:r><w:t>Además el uso de variables slack nos sirve en casos en los que l</w:t></w:r><w:r><w:t xml:space="preserve">os datos del mundo real </w:t></w:r><w:r><w:t>están desordenados</w:t></w:r><w:r><w:t xml:space="preserve"> y casi siempre habrá algún caso que el clasificador no puede acertar, puesto que realizar una separación perfecto no siempre se posible, y en el supuesto de que lo sea, el resultado del modelo puede no ser</w:t></w:r><w:r><w:t xml:space="preserve"> generalizado por otros datos (overfitting). Por lo que p</w:t></w:r><w:r><w:t>ara solucionar este problema y per</w:t></w:r><w:r><w:t>mitir cierta flexibilidad las svm</w:t></w:r><w:r><w:t xml:space="preserve"> utilizan un parámetro C que controla la compensación entre errores de entrenamiento y los márgenes rígidos creando así un </w:t></w:r><w:r><w:t xml:space="preserve">soft-margin </w:t></w:r><w:r><w:t xml:space="preserve">que permite algunos </w:t></w:r><w:r><w:t>errores</w:t></w:r><w:r><w:t xml:space="preserve"> en la clasificación a la </w:t></w:r><w:r><w:t>vez</w:t></w:r><w:r><w:t xml:space="preserve"> que los penaliza.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Cuando esta C </w:t></w:r><w:r><w:t>e</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> pequeña, los errores de clasificación tienen menos importancia y </w:t></w:r><w:r><w:t>el enfoque se encuentra en</w:t></w:r><w:r><w:t xml:space="preserve"> maximizar el margen, mientra</w:t></w:r><w:r><w:t xml:space="preserve">s que cuando C </w:t></w:r><w:r><w:t>e</w:t></w:r><w:r><w:t>s grande, el enfoque se encuentra en</w:t></w:r><w:r><w:t xml:space="preserve"> evitar la clasificación errónea a expensas de mantener el margen pequeño.</w:t></w:r></w:p><w:p><w:r><w:t>Trata de un</w:t></w:r><w:r><w:t xml:space="preserve"> compromiso, o</w:t></w:r><w:r><w:t>b</w:t></w:r><w:r><w:t>tener un mejor clasificador y más</w:t></w:r><w:r><w:t xml:space="preserve"> robusto a expensas de un margen amplio</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Conociendo esto se aplica el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>svc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> con kernel lineal, con </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rbf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>polinomial</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> con grado 3 sobre el conjunto total de datos teniendo en cuenta las mismas variables que en los modelos anteriores.</w:t></w:r></w:p>'
$insertionPoint.InsertXML($newParagraphsXml)
